# The deck's slide-master theme ("ppt/theme/theme1.xml") currently holds the
# "Integral" color scheme. The commit swaps it for the stock PowerPoint
# "Office Theme" color scheme (font scheme and format/effect scheme were
# already byte-identical between the two themes, so only the 12 color-scheme
# slots actually change).
#
# PowerPoint exposes those 12 slots as Design.SlideMaster.Theme.ThemeColorScheme.
# Colors(1..12), in the fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#   9 accent5, 10 accent6, 11 hlink, 12 folHlink
# .RGB takes a packed BGR integer (R + G*256 + B*65536), matching VBA's RGB().

$p = $ppt.ActivePresentation

$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Target "Office Theme" palette (hex RRGGBB), packed as BGR ints for .RGB.
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
